$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(97, 4).Value = 44705

$ws.Cells.Item(98, 4).Value = 44705

$ws.Cells.Item(99, 4).Value = 44567

$ws.Cells.Item(100, 4).Value = 44567

$ws.Cells.Item(101, 4).Value = 44322

$ws.Cells.Item(102, 4).Value = 44322

$ws.Cells.Item(103, 4).Value = 44327

$ws.Cells.Item(104, 4).Value = 44327

$ws.Cells.Item(105, 4).Value = 44383

$ws.Cells.Item(106, 4).Value = 44383

$ws.Cells.Item(107, 4).Value = 44362

$ws.Cells.Item(108, 4).Value = 44362

$ws.Cells.Item(109, 4).Value = 44266

$ws.Cells.Item(110, 4).Value = 44266

$ws.Cells.Item(111, 4).Value = 44607

$ws.Cells.Item(112, 4).Value = 44607

$ws.Cells.Item(113, 4).Value = 44237

$ws.Cells.Item(114, 4).Value = 44237

$ws.Cells.Item(115, 4).Value = 44252

$ws.Cells.Item(116, 4).Value = 44252

$ws.Cells.Item(117, 4).Value = 44453

$ws.Cells.Item(118, 4).Value = 44453

$ws.Cells.Item(119, 4).Value = 44616

$ws.Cells.Item(120, 4).Value = 44616

$ws.Cells.Item(121, 4).Value = 44341
$ws.Cells.Item(121, 10).Value = 200
$ws.Cells.Item(121, 11).Value = 600
$ws.Cells.Item(121, 12).Value = 700
$ws.Cells.Item(121, 13).Value = 650
$ws.Cells.Item(121, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(121, 15).Value = "Región de Ñuble"
$ws.Cells.Item(121, 16).Value = 650
$ws.Cells.Item(121, 17).Value = 1

$ws.Cells.Item(122, 4).Value = 44341
$ws.Cells.Item(122, 9).Value = "Segunda"
$ws.Cells.Item(122, 10).Value = 100
$ws.Cells.Item(122, 11).Value = 500
$ws.Cells.Item(122, 12).Value = 500
$ws.Cells.Item(122, 13).Value = 500
$ws.Cells.Item(122, 16).Value = 500

$ws.Cells.Item(123, 4).Value = 44685
$ws.Cells.Item(123, 9).Value = "Primera"
$ws.Cells.Item(123, 10).Value = 170
$ws.Cells.Item(123, 11).Value = 5500
$ws.Cells.Item(123, 12).Value = 6000
$ws.Cells.Item(123, 13).Value = 5765
$ws.Cells.Item(123, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(123, 15).Value = "Región Metropolitana"
$ws.Cells.Item(123, 16).Value = 160
$ws.Cells.Item(123, 17).Value = 36

$ws.Cells.Item(124, 4).Value = 44609

$ws.Cells.Item(125, 4).Value = 44609

$ws.Cells.Item(126, 4).Value = 44330

$ws.Cells.Item(127, 4).Value = 44330

$ws.Cells.Item(128, 4).Value = 44250

$ws.Cells.Item(129, 4).Value = 44250

$ws.Cells.Item(130, 4).Value = 44334

$ws.Cells.Item(131, 4).Value = 44334

$ws.Cells.Item(132, 4).Value = 44280

$ws.Cells.Item(133, 4).Value = 44280

$ws.Cells.Item(134, 4).Value = 44582

$ws.Cells.Item(135, 4).Value = 44582

$ws.Cells.Item(136, 4).Value = 44257

$ws.Cells.Item(137, 4).Value = 44257

$ws.Cells.Item(138, 4).Value = 44209

$ws.Cells.Item(139, 4).Value = 44209

$ws.Cells.Item(140, 4).Value = 44217

$ws.Cells.Item(141, 4).Value = 44217

$ws.Cells.Item(142, 4).Value = 44405

$ws.Cells.Item(143, 4).Value = 44405

$ws.Cells.Item(144, 4).Value = 44475

$ws.Cells.Item(145, 4).Value = 44475

$ws.Cells.Item(146, 4).Value = 44239

$ws.Cells.Item(147, 4).Value = 44239

$ws.Cells.Item(148, 4).Value = 44358

$ws.Cells.Item(149, 4).Value = 44358

$ws.Cells.Item(150, 4).Value = 44187

$ws.Cells.Item(151, 4).Value = 44187

$ws.Cells.Item(152, 4).Value = 44694

$ws.Cells.Item(153, 4).Value = 44694

$ws.Cells.Item(154, 4).Value = 44488

$ws.Cells.Item(155, 4).Value = 44488

$ws.Cells.Item(156, 4).Value = 44166
$ws.Cells.Item(156, 15).Value = "Región de Ñuble"

$ws.Cells.Item(157, 4).Value = 44166
$ws.Cells.Item(157, 15).Value = "Región de Ñuble"

$ws.Cells.Item(158, 4).Value = 44316
$ws.Cells.Item(158, 10).Value = 200
$ws.Cells.Item(158, 11).Value = 600
$ws.Cells.Item(158, 12).Value = 700
$ws.Cells.Item(158, 13).Value = 650
$ws.Cells.Item(158, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(158, 16).Value = 650
$ws.Cells.Item(158, 17).Value = 1

$ws.Cells.Item(159, 4).Value = 44316
$ws.Cells.Item(159, 9).Value = "Segunda"
$ws.Cells.Item(159, 10).Value = 100
$ws.Cells.Item(159, 11).Value = 500
$ws.Cells.Item(159, 12).Value = 500
$ws.Cells.Item(159, 13).Value = 500
$ws.Cells.Item(159, 15).Value = "Región Metropolitana"
$ws.Cells.Item(159, 16).Value = 500

$ws.Cells.Item(160, 4).Value = 44656
$ws.Cells.Item(160, 9).Value = "Primera"
$ws.Cells.Item(160, 10).Value = 130
$ws.Cells.Item(160, 11).Value = 650
$ws.Cells.Item(160, 12).Value = 6000
$ws.Cells.Item(160, 13).Value = 3942
$ws.Cells.Item(160, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(160, 15).Value = "Región Metropolitana"
$ws.Cells.Item(160, 16).Value = 110
$ws.Cells.Item(160, 17).Value = 36

$ws.Cells.Item(161, 4).Value = 44469
$ws.Cells.Item(161, 15).Value = "Región de Ñuble"

$ws.Cells.Item(162, 4).Value = 44469
$ws.Cells.Item(162, 15).Value = "Región de Ñuble"

$ws.Cells.Item(163, 4).Value = 44579
$ws.Cells.Item(163, 15).Value = "Región Metropolitana"

$ws.Cells.Item(164, 4).Value = 44579
$ws.Cells.Item(164, 15).Value = "Región Metropolitana"

$ws.Cells.Item(165, 4).Value = 44267
$ws.Cells.Item(165, 10).Value = 200
$ws.Cells.Item(165, 11).Value = 600
$ws.Cells.Item(165, 12).Value = 700
$ws.Cells.Item(165, 13).Value = 650
$ws.Cells.Item(165, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(165, 15).Value = "Región de Ñuble"
$ws.Cells.Item(165, 16).Value = 650
$ws.Cells.Item(165, 17).Value = 1

$ws.Cells.Item(166, 4).Value = 44267
$ws.Cells.Item(166, 9).Value = "Segunda"
$ws.Cells.Item(166, 10).Value = 100
$ws.Cells.Item(166, 11).Value = 500
$ws.Cells.Item(166, 12).Value = 500
$ws.Cells.Item(166, 13).Value = 500
$ws.Cells.Item(166, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(166, 15).Value = "Región de Ñuble"
$ws.Cells.Item(166, 16).Value = 500
$ws.Cells.Item(166, 17).Value = 1

$ws.Cells.Item(167, 4).Value = 44671
$ws.Cells.Item(167, 10).Value = 110
$ws.Cells.Item(167, 11).Value = 6500
$ws.Cells.Item(167, 12).Value = 7000
$ws.Cells.Item(167, 13).Value = 6773
$ws.Cells.Item(167, 14).Value = "$/docena de atados"
$ws.Cells.Item(167, 15).Value = "Región Metropolitana"
$ws.Cells.Item(167, 16).Value = 2258
$ws.Cells.Item(167, 17).Value = 3

$ws.Cells.Item(168, 4).Value = 44672
$ws.Cells.Item(168, 9).Value = "Primera"
$ws.Cells.Item(168, 10).Value = 150
$ws.Cells.Item(168, 11).Value = 4500
$ws.Cells.Item(168, 12).Value = 5000
$ws.Cells.Item(168, 13).Value = 4767
$ws.Cells.Item(168, 14).Value = "$/caja 36 atados"
$ws.Cells.Item(168, 15).Value = "Región Metropolitana"
$ws.Cells.Item(168, 16).Value = 132
$ws.Cells.Item(168, 17).Value = 36

$ws.Cells.Item(169, 4).Value = 44490

$ws.Cells.Item(170, 4).Value = 44490

$ws.Cells.Item(171, 4).Value = 44679

$ws.Cells.Item(172, 4).Value = 44679

$ws.Cells.Item(173, 4).Value = 44365

$ws.Cells.Item(174, 4).Value = 44365

$ws.Cells.Item(175, 4).Value = 44427

$ws.Cells.Item(176, 4).Value = 44427

$ws.Cells.Item(177, 4).Value = 44565

$ws.Cells.Item(178, 4).Value = 44565

$ws.Cells.Item(179, 4).Value = 44447

$ws.Cells.Item(180, 4).Value = 44447

$ws.Cells.Item(181, 4).Value = 44523

$ws.Cells.Item(182, 4).Value = 44523

$ws.Cells.Item(183, 4).Value = 44462

$ws.Cells.Item(184, 4).Value = 44462

$ws.Cells.Item(185, 4).Value = 44159

$ws.Cells.Item(186, 4).Value = 44159

$ws.Cells.Item(187, 4).Value = 44344

$ws.Cells.Item(188, 4).Value = 44344

$ws.Cells.Item(189, 4).Value = 44692

$ws.Cells.Item(190, 4).Value = 44692

$ws.Cells.Item(191, 4).Value = 44376

$ws.Cells.Item(192, 4).Value = 44376

$ws.Cells.Item(193, 1).Value = 11
$ws.Cells.Item(193, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(193, 3).Value = "Bíobío"
$ws.Cells.Item(193, 4).Value = 44442
$ws.Cells.Item(193, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(193, 5).Value = 8
$ws.Cells.Item(193, 6).Value = 100112040
$ws.Cells.Item(193, 7).Value = "Cilantro"
$ws.Cells.Item(193, 8).Value = "Sin especificar"
$ws.Cells.Item(193, 9).Value = "Primera"
$ws.Cells.Item(193, 10).Value = 200
$ws.Cells.Item(193, 11).Value = 600
$ws.Cells.Item(193, 12).Value = 700
$ws.Cells.Item(193, 13).Value = 650
$ws.Cells.Item(193, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(193, 15).Value = "Región de Ñuble"
$ws.Cells.Item(193, 16).Value = 650
$ws.Cells.Item(193, 17).Value = 1
$ws.Cells.Item(193, 18).Value = "Hortaliza"

$ws.Cells.Item(194, 1).Value = 11
$ws.Cells.Item(194, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(194, 3).Value = "Bíobío"
$ws.Cells.Item(194, 4).Value = 44442
$ws.Cells.Item(194, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(194, 5).Value = 8
$ws.Cells.Item(194, 6).Value = 100112040
$ws.Cells.Item(194, 7).Value = "Cilantro"
$ws.Cells.Item(194, 8).Value = "Sin especificar"
$ws.Cells.Item(194, 9).Value = "Segunda"
$ws.Cells.Item(194, 10).Value = 100
$ws.Cells.Item(194, 11).Value = 500
$ws.Cells.Item(194, 12).Value = 500
$ws.Cells.Item(194, 13).Value = 500
$ws.Cells.Item(194, 14).Value = "$/atado 0,5 a 1 kilo"
$ws.Cells.Item(194, 15).Value = "Región de Ñuble"
$ws.Cells.Item(194, 16).Value = 500
$ws.Cells.Item(194, 17).Value = 1
$ws.Cells.Item(194, 18).Value = "Hortaliza"
